$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old row 17 ("eth" / Significant Related Persons / "Vitalik Buterin")
# so every subsequent "btc" row shifts up by one.
$ws.Rows.Item(17).Delete()

# Append the new "eth" related-persons rows at the bottom (now rows 24-28).
$ws.Range("A24").Value = "eth"
$ws.Range("B24").Value = "Significant Related Persons"
$ws.Range("C24").Value = "VitalikButerin"
$ws.Range("D24").Value = "test111"

$ws.Range("A25").Value = "eth"
$ws.Range("B25").Value = "Significant Related Persons"
$ws.Range("C25").Value = "ethereum"
$ws.Range("D25").Value = "test111"

$ws.Range("A26").Value = "eth"
$ws.Range("B26").Value = "Significant Related Persons"
$ws.Range("C26").Value = "brian_armstrong"
$ws.Range("D26").Value = "test111"

$ws.Range("A27").Value = "eth"
$ws.Range("B27").Value = "Significant Related Persons"
$ws.Range("C27").Value = "CoinDesk"
$ws.Range("D27").Value = "test111"

$ws.Range("A28").Value = "eth"
$ws.Range("B28").Value = "Significant Related Persons"
$ws.Range("C28").Value = "vip1"
$ws.Range("D28").Value = "test111"
